$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 0.03375954870909267, 0.6485414634705653, 0.7314436294522902, 0.8552447775065862, 0.8663660556530268, 37),
    @(3, 0.03135085916488722, 0.5511937747325102, 0.5928526553971856, 0.7699692561376626, 0.7802437722544544, 36),
    @(4, 0.0324755997945022,  0.5741707009966609, 0.6201381699606496, 0.7874885205262675, 0.7983056080049333, 35),
    @(5, 0.04451789004105118, 0.5528601565202774, 0.615403550752249,  0.7844766094360296, 0.7949907199654754, 34),
    @(6, 0.02488288622478917, 0.5597336268738891, 0.6071832843101319, 0.7792196637085924, 0.7908977527023375, 33),
    @(7, 0.06447809342936392, 0.5639312790710509, 0.6386050102870782, 0.7991276558141873, 0.8092673457567485, 32),
    @(8, 0.0316477260114686,  0.5875569779702665, 0.6486140808091702, 0.8053658055872314, 0.8180462001963719, 31),
    @(9, 0.05575202370392403, 0.5712956843590213, 0.6777300950868858, 0.8232436426033825, 0.8353948933700528, 30),
    @(10,0.03351674390117437, 0.5951644878936757, 0.6821752870341646, 0.8259390334850174, 0.8398661783807415, 29),
    @(11,0.06637440518007821, 0.5683437572335094, 0.6866267435578414, 0.8286294368159035, 0.8411234351010097, 28)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}
